$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Sin ítem - Sub. 27"
$ws.Range("H2").Value = "Sin ítem - Sub. 28"
$ws.Range("O2").Value = "Sin ítem - Sub. 35"
